# Observations-summary sheet: drop the "documentation2" and "presence" demo
# rows, keeping only the (fixed-up) "documentation" profile row.
#
# Original data rows (row 1 is the header):
#   row2: us-core-observation-adi-documentation  ... Code="" / Code VS=<LOINC valueset URL>
#   row3: us-core-observation-adi-documentation2 ... Code="LOINC#45473-6" / Code VS=""
#   row4: us-core-observation-adi-presence       ... Code="LOINC#45473-6" / Code VS=""
#
# Target: a single data row equal to row 3's Code/Code VS pairing (Code =
# "LOINC#45473-6", Code VS = "") but row 2's Profile/Name labels. We delete
# row 2 (the stale Code/Code VS pairing) and row 4 (the extra "presence"
# row), which shifts the still-correct row 3 up into row 2, then fix up
# just the Profile/Name columns on the surviving row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 2 (shifts row 3 -> row 2, row 4 -> row 3).
$ws.Range("A2:K2").Delete()
# Remove the old row 4, now sitting at row 3 (shifts nothing else up, it's
# the last row).
$ws.Range("A3:K3").Delete()

# The surviving row (originally row 3) already carries the right Code /
# Code VS values; just correct the Profile id and display Name back to the
# non-"2" variant.
$ws.Range("A2").Value2 = "us-core-observation-adi-documentation"
$ws.Range("B2").Value2 = "US Core Observation ADI Documentation Profile"
